# Apply the "Add files via upload" revision to the ProliferationLevel sheet:
# the row group for code 713 (rows 1219-1228, one row per year 1970-1979) was
# removed, so every following row shifts up by 10 (1496 rows -> 1486 rows).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ProliferationLevel")

# Remove the 10-row block for code 713; everything below shifts up to fill
# the gap (codes/years/values are otherwise untouched).
$ws.Range("A1219:C1228").EntireRow.Delete()

# The AutoFilter range shrinks along with the data (A1:C1496 -> A1:C1486).
$ws.AutoFilterMode = $false
$ws.Range("A1:C1486").AutoFilter()

# Keep the workbook-level hidden _FilterDatabase name for this sheet in sync.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
  $n = $wb.Names.Item($i)
  if ($n.Name -eq "ProliferationLevel!_FilterDatabase") {
    $n.RefersTo = "=ProliferationLevel!`$A`$1:`$C`$1486"
  }
}

# ProliferationLevel becomes the active sheet/tab (was TreatyKey), with the
# selection left on D1229 (the row that used to be 731/1970 before shifting).
$ws.Activate()
$ws.Range("D1229").Select()
